# Week 13 logging update
$wb = $excel.ActiveWorkbook

# --- YDS sheet: append per-game yard figures to the running season lists ---
$ws = $wb.Worksheets.Item("YDS")
$ws.Range("B2").Value = $ws.Range("B2").Value() + " 0 5 9 2 2 1 4 9 7 18 9 4 1 3 0 7 6 2 2 6 0 0 1 1 3 -1 3 5 -1 8 34 1 4 2 2 7 12 3 2 0"
$ws.Range("C2").Value = $ws.Range("C2").Value() + " 2 1 3 17 9 1 0 1 15 2 2 3 2 3 8 17"
$ws.Range("B3").Value = $ws.Range("B3").Value() + " 16 9 36 12 10 12 5 25 13 16 12 22 -2 2 28 9 7 7 2 1"
$ws.Range("C3").Value = $ws.Range("C3").Value() + " 8 7 3 10 7 5 10 29 4 19 5 18 14 19 11 -2 15 4 14 7 11"

# --- OFF sheet: totals for Home / Road rows ---
$ws = $wb.Worksheets.Item("OFF")
$ws.Range("C2").Value = 170
$ws.Range("D2").Value = 12
$ws.Range("E2").Value = 9
$ws.Range("F2").Value = 59
$ws.Range("G2").Value = 54
$ws.Range("I2").Value = 11
$ws.Range("J2").Value = 36
$ws.Range("N2").Value = 11
$ws.Range("O2").Value = 15
$ws.Range("P2").Value = 5

$ws.Range("C3").Value = 113
$ws.Range("E3").Value = 23
$ws.Range("F3").Value = 69
$ws.Range("G3").Value = 19
$ws.Range("H3").Value = 16
$ws.Range("I3").Value = 30
$ws.Range("J3").Value = 44
$ws.Range("L3").Value = 236
$ws.Range("M3").Value = 151
$ws.Range("Q3").Value = 487

# --- DEF sheet: totals for Home / Road rows ---
$ws = $wb.Worksheets.Item("DEF")
$ws.Range("C2").Value = 162
$ws.Range("D2").Value = 10
$ws.Range("F2").Value = 41
$ws.Range("G2").Value = 56
$ws.Range("N2").Value = 19
$ws.Range("O2").Value = 22
$ws.Range("P2").Value = 10

$ws.Range("C3").Value = 145
$ws.Range("E3").Value = 26
$ws.Range("F3").Value = 93
$ws.Range("I3").Value = 50
$ws.Range("J3").Value = 38
$ws.Range("L3").Value = 269
$ws.Range("M3").Value = 183
$ws.Range("Q3").Value = 486

# --- ST sheet ---
$ws = $wb.Worksheets.Item("ST")
$ws.Range("B2").Value = 69
$ws.Range("D2").Value = 43
$ws.Range("F2").Value = 173
$ws.Range("G2").Value = 164
$ws.Range("J2").Value = 66
$ws.Range("K2").Value = 60
$ws.Range("L2").Value = 42
$ws.Range("M2").Value = 36
$ws.Range("B3").Value = 45
$ws.Range("B4").Value = $ws.Range("B4").Value() + " 65 61 67"
$ws.Range("B5").Value = $ws.Range("B5").Value() + " 79 28 32"
$ws.Range("B6").Value = $ws.Range("B6").Value() + " 21 5 13"
$ws.Range("D3").Value = $ws.Range("D3").Value() + " 38"
$ws.Range("D4").Value = $ws.Range("D4").Value() + " 0"
$ws.Range("D5").Value = $ws.Range("D5").Value() + " 0"

# --- TURNS sheet ---
$ws = $wb.Worksheets.Item("TURNS")
$ws.Range("C3").Value = 5
$ws.Range("D3").Value = 6

# --- PEN sheet ---
$ws = $wb.Worksheets.Item("PEN")
$ws.Range("D2").Value = 10
$ws.Range("B3").Value = 10
$ws.Range("D3").Value = 5
